# Check_ThatUserCanSearchByEdited_GF_Name.xlsx - Q2 search test data update.
# Renames the "symbol" / "name" / "editedname" sample values on row 2 and
# widens columns E/F slightly, per the new GF_Tag/name_GF_edite* test case.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# Row 2 sample values (symbol / name / editedname columns).
$ws.Range("E2").Value = "GF_Tag_edite_selenium"
$ws.Range("F2").Value = "name_GF_edite"
$ws.Range("H2").Value = "name_GF_edited"

# Column widths: target stored (OOXML) widths are 24 and 19.28515625 chars.
# ColumnWidth here is quantized to 1/6-character steps, so feed it the value
# whose quantization lands exactly on the desired stored width (24 -> 139/6).
$ws.Columns.Item(5).ColumnWidth = 23.166666666666668
$ws.Columns.Item(6).ColumnWidth = 19.28515625

# Move the active selection to G7.
$ws.Range("G7").Select()
